$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new patient ("Vier") record needs to be inserted as row 5, pushing the
# existing last patient record (Simon Wosnitza) down to row 6.
#
# Writing directly into row 5 (which already holds data) would keep that
# row's existing cell formatting; writing into row 6 (which is currently
# blank / past the used range) lets the new cells pick up the worksheet's
# column default formatting - which is what the new row ends up using.
# So: write the new record into row 6 first, stash it in a scratch row,
# move the old row 5 down into row 6, then drop the new record into row 5.

$ws.Range("A6").Value = "Fgk4hj"
$ws.Range("B6").Value = "hDfzu85Rf7"
$ws.Range("C6").Value = "Patient"
$ws.Range("D6").Value = "Vier"
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = "07/01/2022"
$ws.Range("H6").Value = "07/04/2022"

$ws.Range("A6:H6").Copy($ws.Range("A7:H7"))
$ws.Range("A5:H5").Copy($ws.Range("A6:H6"))
$ws.Range("A7:H7").Copy($ws.Range("A5:H5"))
$ws.Range("A7:H7").Delete()

# Restore the row heights used throughout the table.
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item(2).RowHeight = 18.75
$ws.Rows.Item(3).RowHeight = 18.75
$ws.Rows.Item(4).RowHeight = 18.75
$ws.Rows.Item(5).RowHeight = 18
$ws.Rows.Item(6).RowHeight = 18.75
